$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "/flashcard.html" value into G9 (reuses existing shared string)
$ws.Range("G9").Value = "/flashcard.html"

# Update the active selection to F14, matching the saved cursor position
$ws.Range("F14").Select() | Out-Null
